# Generate Report for Handback
# Updates the localization-status report after a handback run:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet (zh-cn / de-de columns) and on each language sheet.
#  - The "Latest Handback DateTime" for each language is refreshed.
#  - The stale "handback file is not latest" Error Detail warning is cleared
#    now that the handback is in sync.
#  - The Status / Error Detail columns are widened to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$overview.Range("E2").Value2 = $newStatus
$overview.Range("F2").Value2 = $newStatus

# Widen the Status columns (E, F) to fit the longer text
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# --- zh-cn sheet ------------------------------------------------------
$zhcn.Range("C2").Value2 = $newStatus
$zhcn.Range("K2").Value2 = "2016-08-27 10:49:14"
$zhcn.Range("P2").Value2 = ""

$zhcn.Columns.Item(3).ColumnWidth  = 29.1
$zhcn.Columns.Item(16).ColumnWidth = 12.76

# --- de-de sheet ------------------------------------------------------
$dede.Range("C2").Value2 = $newStatus
$dede.Range("K2").Value2 = "2016-08-27 10:49:21"
$dede.Range("P2").Value2 = ""

$dede.Columns.Item(3).ColumnWidth  = 29.1
$dede.Columns.Item(16).ColumnWidth = 12.76
